$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.333.57"
$ws.Range("D2").Style = $ws.Range("B2").Style

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.07%  "
$ws.Range("E2").Style = $ws.Range("B2").Style

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.773.77"
$ws.Range("D3").Style = $ws.Range("B3").Style

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("E3").Style = $ws.Range("B3").Style

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = $ws.Range("B4").Style

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E4").Style = $ws.Range("B4").Style

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E5").Style = $ws.Range("B5").Style

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.76"
$ws.Range("D6").Style = $ws.Range("B6").Style

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E6").Style = $ws.Range("B6").Style

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4228"
$ws.Range("D7").Style = $ws.Range("B7").Style

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("E7").Style = $ws.Range("B7").Style

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("E8").Style = $ws.Range("B8").Style

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07129"
$ws.Range("D9").Style = $ws.Range("B9").Style

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("E9").Style = $ws.Range("B9").Style

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8357"
$ws.Range("D10").Style = $ws.Range("B10").Style

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("E10").Style = $ws.Range("B10").Style

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.38"
$ws.Range("D11").Style = $ws.Range("B11").Style

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("E11").Style = $ws.Range("B11").Style

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.762.76"
$ws.Range("D12").Style = $ws.Range("B12").Style

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.57%  "
$ws.Range("E12").Style = $ws.Range("B12").Style

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.446"
$ws.Range("D13").Style = $ws.Range("B13").Style

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("E13").Style = $ws.Range("B13").Style

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.236"
$ws.Range("D14").Style = $ws.Range("B14").Style

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("E14").Style = $ws.Range("B14").Style

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06867"
$ws.Range("D15").Style = $ws.Range("B15").Style

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("E15").Style = $ws.Range("B15").Style

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = $ws.Range("B16").Style

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("E16").Style = $ws.Range("B16").Style

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "78.84"
$ws.Range("D17").Style = $ws.Range("B17").Style

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.02%  "
$ws.Range("E17").Style = $ws.Range("B17").Style

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008633"
$ws.Range("D18").Style = $ws.Range("B18").Style

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E18").Style = $ws.Range("B18").Style

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("E19").Style = $ws.Range("B19").Style

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.91"
$ws.Range("D20").Style = $ws.Range("B20").Style

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("E20").Style = $ws.Range("B20").Style

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.347.83"
$ws.Range("D21").Style = $ws.Range("B21").Style

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.57%  "
$ws.Range("E21").Style = $ws.Range("B21").Style

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.080"
$ws.Range("D22").Style = $ws.Range("B22").Style

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("E22").Style = $ws.Range("B22").Style

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.92"
$ws.Range("D23").Style = $ws.Range("B23").Style

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("E23").Style = $ws.Range("B23").Style

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.005.34"
$ws.Range("D24").Style = $ws.Range("B24").Style

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("E24").Style = $ws.Range("B24").Style

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.47"
$ws.Range("D25").Style = $ws.Range("B25").Style

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("E25").Style = $ws.Range("B25").Style

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -8.76%  "
$ws.Range("E26").Style = $ws.Range("B26").Style

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.03"
$ws.Range("D27").Style = $ws.Range("B27").Style

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E27").Style = $ws.Range("B27").Style

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.055"
$ws.Range("D28").Style = $ws.Range("B28").Style

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E28").Style = $ws.Range("B28").Style

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.32"
$ws.Range("D29").Style = $ws.Range("B29").Style

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("E29").Style = $ws.Range("B29").Style

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.828"
$ws.Range("D30").Style = $ws.Range("B30").Style

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +9.38%  "
$ws.Range("E30").Style = $ws.Range("B30").Style

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08827"
$ws.Range("D31").Style = $ws.Range("B31").Style

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E31").Style = $ws.Range("B31").Style

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7268"
$ws.Range("D32").Style = $ws.Range("B32").Style

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("E32").Style = $ws.Range("B32").Style

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.114"
$ws.Range("D33").Style = $ws.Range("B33").Style

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("E33").Style = $ws.Range("B33").Style

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.320"
$ws.Range("D34").Style = $ws.Range("B34").Style

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("E34").Style = $ws.Range("B34").Style

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.001"
$ws.Range("D35").Style = $ws.Range("B35").Style

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E35").Style = $ws.Range("B35").Style

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.735"
$ws.Range("D36").Style = $ws.Range("B36").Style

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.27%  "
$ws.Range("E36").Style = $ws.Range("B36").Style

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.085"
$ws.Range("D37").Style = $ws.Range("B37").Style

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("E37").Style = $ws.Range("B37").Style

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05134"
$ws.Range("D38").Style = $ws.Range("B38").Style

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E38").Style = $ws.Range("B38").Style

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01882"
$ws.Range("D39").Style = $ws.Range("B39").Style

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("E39").Style = $ws.Range("B39").Style

$ws.Range("B40").Value = "TheSandbox"

$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4909"
$ws.Range("D40").Style = $ws.Range("B40").Style

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("E40").Style = $ws.Range("B40").Style

$ws.Range("B41").Value = "Algorand"

$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1607"
$ws.Range("D41").Style = $ws.Range("B41").Style

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("E41").Style = $ws.Range("B41").Style

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.598"
$ws.Range("D42").Style = $ws.Range("B42").Style

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.73%  "
$ws.Range("E42").Style = $ws.Range("B42").Style

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.361"
$ws.Range("D43").Style = $ws.Range("B43").Style

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.99%  "
$ws.Range("E43").Style = $ws.Range("B43").Style

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.947"
$ws.Range("D44").Style = $ws.Range("B44").Style

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.92%  "
$ws.Range("E44").Style = $ws.Range("B44").Style

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.61"
$ws.Range("D45").Style = $ws.Range("B45").Style

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("E45").Style = $ws.Range("B45").Style

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = $ws.Range("B46").Style

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E46").Style = $ws.Range("B46").Style

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.15"
$ws.Range("D47").Style = $ws.Range("B47").Style

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("E47").Style = $ws.Range("B47").Style

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.630"
$ws.Range("D48").Style = $ws.Range("B48").Style

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.41%  "
$ws.Range("E48").Style = $ws.Range("B48").Style

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06175"
$ws.Range("D49").Style = $ws.Range("B49").Style

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("E49").Style = $ws.Range("B49").Style

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4456"
$ws.Range("D50").Style = $ws.Range("B50").Style

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("E50").Style = $ws.Range("B50").Style

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.717"
$ws.Range("D51").Style = $ws.Range("B51").Style

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.93%  "
$ws.Range("E51").Style = $ws.Range("B51").Style
